# CRUD operations on Teams
# Update planned/realized hours for several tasks in the "Planning & Journal"
# sheet, highlight the two corrected "realised hours" cells in red, and move
# the current selection to reflect where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71 - "[S3] Intégrer les différents endpoints de l'API aux interfaces"
$ws.Range("D71").Value = 4
$ws.Range("E71").Value = 3.5

# Row 72 - "Amélioration du chargement initial des données"
$ws.Range("D72").Value = 4

# Row 73 - "[S3] Mettre en place la connexion entre les interface et l'API"
# (realised hours corrected and highlighted in red)
$ws.Range("E73").Value = 4
$ws.Range("E73").Font.Color = 255

# Row 74 - "[S3] Intégrer les différents endpoints de l'API aux interfaces"
# (planned + realised hours corrected, realised hours highlighted in red)
$ws.Range("D74").Value = 24
$ws.Range("E74").Value = 12
$ws.Range("E74").Font.Color = 255

# Row 76 - realised hours corrected
$ws.Range("E76").Value = 20

# Row 77 - realised hours corrected
$ws.Range("E77").Value = 40

# Update the current selection / scroll position to match where the author
# ended up after making these edits.
[void]$ws.Activate()
[void]$ws.Range("H76").Select()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 1
